# Update the title placeholder in cell A1 (merged A1:G1) of the export
# template: replace the static Vietnamese heading text with a template
# placeholder "{{Name}}".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "{{Name}}"

# Update the active selection to match the merged header range (A1:G1)
# instead of the previously selected cell F3.
$ws.Range("A1:G1").Select()
